$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove 3 worker/period rows from the data table (table shrinks from 11 to 8 rows)
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(17).Delete()

# Update summary figures
$ws.Range("E11").Value = 128748
$ws.Range("C13").Value = 5
$ws.Range("F13").Value = 7

# Rewrite the worker/period detail rows with the new data set
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1051450548"
$ws.Range("D16").Value = "YULISA DIAZ CUADRADO"
$ws.Range("E16").Value = "2110"
$ws.Range("F16").Value = 30400
$ws.Range("G16").Value = 1200000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1051450548"
$ws.Range("D17").Value = "YULISA DIAZ CUADRADO"
$ws.Range("E17").Value = "2111"
$ws.Range("F17").Value = 48000
$ws.Range("G17").Value = 1200000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1051450548"
$ws.Range("D18").Value = "YULISA DIAZ CUADRADO"
$ws.Range("E18").Value = "2112"
$ws.Range("F18").Value = 9600
$ws.Range("G18").Value = 1200000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "45503050"
$ws.Range("D19").Value = "MERCEDES CARO CARO"
$ws.Range("E19").Value = "2204"
$ws.Range("F19").Value = 5867
$ws.Range("G19").Value = 1100000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "73125000"
$ws.Range("D20").Value = "FREDY DE JESUS GUARDO GAVIRIA"
$ws.Range("E20").Value = "2205"
$ws.Range("F20").Value = 14667
$ws.Range("G20").Value = 1000000

$ws.Range("B21").Value = "PPT"
$ws.Range("C21").Value = "5484560"
$ws.Range("D21").Value = "JHOANA ALICIA PEROZO URDANTEA"
$ws.Range("E21").Value = "2205"
$ws.Range("F21").Value = 14667
$ws.Range("G21").Value = 1000000

$ws.Range("B22").Value = "PPT"
$ws.Range("C22").Value = "5484560"
$ws.Range("D22").Value = "JHOANA ALICIA PEROZO URDANTEA"
$ws.Range("E22").Value = "2206"
$ws.Range("F22").Value = 4000
$ws.Range("G22").Value = 1000000

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "73099133"
$ws.Range("D23").Value = "WILLIAM LARA MARIN"
$ws.Range("E23").Value = "2303"
$ws.Range("F23").Value = 1547
$ws.Range("G23").Value = 1160000

# Column D auto-fits to the now-shorter longest worker name
$ws.Columns.Item(4).ColumnWidth = 31.98
